# Fixed spelling in 5510 video06
#
# Applies the spelling corrections from the commit to the speaker notes of
# several slides, plus one on-slide label fix (slide 30: "Nonrandomized"
# -> "Non-randomized").

$p = $ppt.ActivePresentation

function Fix-NotesText($SlideIndex, $Pairs) {
    $slide = $p.Slides.Item($SlideIndex)
    $notesPage = $slide.NotesPage

    $target = $null
    for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
        $shape = $notesPage.Shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText -and $shape.Name -like "Notes Placeholder*") {
            $target = $shape
            break
        }
    }
    if ($null -eq $target) {
        Write-Output "WARNING: no notes placeholder found on slide $SlideIndex"
        return
    }

    $tr = $target.TextFrame.TextRange
    $text = $tr.Text
    for ($i = 0; $i -lt $Pairs.Count; $i += 2) {
        $old = $Pairs[$i]
        $new = $Pairs[$i + 1]
        $text = $text -replace [regex]::Escape($old), $new
    }
    $tr.Text = $text
}

Fix-NotesText 2  @("Althought", "Although")
Fix-NotesText 3  @("strenuosly", "strenuously")
Fix-NotesText 4  @("Helathcare", "Healthcare")
Fix-NotesText 5  @("persepctive.", "perspective.")
Fix-NotesText 6  @("ateempts", "attempts", "processs", "process")
Fix-NotesText 8  @("obective", "objective")
Fix-NotesText 10 @("Oucome", "Outcome", "influece", "influence", "outdome", "outcome")
Fix-NotesText 25 @("soemtimes", "sometimes")
Fix-NotesText 29 @("occured", "occurred")
Fix-NotesText 33 @("occuring", "occurring")
Fix-NotesText 40 @("withdrwal", "withdrawal", "beneifts", "benefits", "plently", "plenty")

# Slide 30: "Content Placeholder 2" body text - fix "Nonrandomized" -> "Non-randomized"
$slide30 = $p.Slides.Item(30)
for ($i = 1; $i -le $slide30.Shapes.Count; $i++) {
    $shape = $slide30.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "*Nonrandomized comparison*") {
            $tr.Text = $tr.Text -replace [regex]::Escape("Nonrandomized comparison"), "Non-randomized comparison"
        }
    }
}
